$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint2")

# ---------------------------------------------------------------------------
# The old row 11 (US20 "Aunts and uncles") is about to be pushed down to
# row 27, with the task breakdown for US19 ("Aunts and uncles should not
# marry their nieces or nephews") inserted above it as rows 11-26, and a new
# task breakdown for US20 inserted below it as rows 28-34.
# ---------------------------------------------------------------------------

# Insert 16 blank rows at row 11 for the T19.xx breakdown (this pushes the
# existing "US20 / Aunts and uncles" row from 11 down to 27 and copies the
# C/D column formatting down automatically).
$ws.Rows("11:26").Insert()

# Insert 7 blank rows after the relocated US20 row (now at row 27) for the
# T20.xx breakdown.
$ws.Rows("28:34").Insert()

# --- Fill in the E/F (Est Size / Est Time) values for the two story rows ---
$ws.Range("E11").Value = 65
$ws.Range("F11").Value = 90
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = 80

# --- Owner / Status columns for every new task row (same as the rest of the sheet) ---
$owner = "zg"
$status = "Assigned"

# --- T19.xx breakdown (rows 11-26) ------------------------------------------------
$ws.Range("A11").Value = "T19.01"
$ws.Range("B11").Value = "Iterate through families"

$ws.Range("A12").Value = "T19.02"
$ws.Range("B12").Value = "Find if there was a marriage"

$ws.Range("A13").Value = "T19.03"
$ws.Range("B13").Value = "Store husband ID"

$ws.Range("A14").Value = "T19.04"
$ws.Range("B14").Value = "Store wife ID"

$ws.Range("A15").Value = "T19.05"
$ws.Range("B15").Value = "Find family ID of husband through child "

$ws.Range("A16").Value = "T19.06"
$ws.Range("B16").Value = "Find family ID of wife through child "

$ws.Range("A17").Value = "T19.07"
$ws.Range("B17").Value = "Find husband's father"

$ws.Range("A18").Value = "T19.08"
$ws.Range("B18").Value = "Find husband's mother"

$ws.Range("A19").Value = "T19.09"
$ws.Range("B19").Value = "Find wife's father"

$ws.Range("A20").Value = "T19.10"
$ws.Range("B20").Value = "Find wife's mother"

$ws.Range("A21").Value = "T19.11"
$ws.Range("B21").Value = "Find family ID of husband's father through child"

$ws.Range("A22").Value = "T19.12"
$ws.Range("B22").Value = "Find family ID of husband's mother through child"

$ws.Range("A23").Value = "T19.13"
$ws.Range("B23").Value = "Find family ID of wife's father through child"

$ws.Range("A24").Value = "T19.14"
$ws.Range("B24").Value = "Find family ID of wife's mother through child"

$ws.Range("A25").Value = "T19.15"
$ws.Range("B25").Value = "Find if family IDs of either parent of husband and wife are the same"

$ws.Range("A26").Value = "T19.16"
$ws.Range("B26").Value = "Output an error if there is an overlap"

# --- T20.xx breakdown (rows 28-34) ------------------------------------------------
$ws.Range("A28").Value = "T20.01"
$ws.Range("B28").Value = "Iterate through family "

$ws.Range("A29").Value = "T20.02"
$ws.Range("B29").Value = "Store children"

$ws.Range("A30").Value = "T20.03"
$ws.Range("B30").Value = "Store husband ID"

$ws.Range("A31").Value = "T20.04"
$ws.Range("B31").Value = "Store wife ID"

$ws.Range("A32").Value = "T20.05"
$ws.Range("B32").Value = "Iterate through individuals to find family ID through child"

$ws.Range("A33").Value = "T20.06"
$ws.Range("B33").Value = "Store children through family ID"

$ws.Range("A34").Value = "T20.07"
$ws.Range("B34").Value = "Store aunt and uncles as children through previous family ID"

# Owner / Status for all new task rows (11-26 and 28-34)
$ws.Range("C11:D26").Value = $owner
$ws.Range("D11:D26").Value = $status
for ($r = 11; $r -le 26; $r++) {
    $ws.Cells.Item($r, 3).Value = $owner
    $ws.Cells.Item($r, 4).Value = $status
}
for ($r = 28; $r -le 34; $r++) {
    $ws.Cells.Item($r, 3).Value = $owner
    $ws.Cells.Item($r, 4).Value = $status
}

# --- Column B is now much wider to fit the longer task descriptions ---
$ws.Columns("B:B").ColumnWidth = 54.61328125

# --- Update the sheet view: scrolled down a bit, zoomed to 100%, and the
#     active selection left on B11 (the first new task row) ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 100
$ws.Range("B11").Select()
